# Working solutions with min of 30 hour per day
#
# 1) Employees sheet: set "min_hours_per_week" (col I) to 30 for every
#    employee row, and remove the last two employees (Jacob, Leah).
# 2) Parameters sheet: set "minusHours_per_Week" (B4) to 0.
# 3) Leave the Employees sheet as the active/selected sheet & tab, with
#    the just-deleted rows highlighted (mirrors the author's recorded
#    selection after deleting the rows).

$wb = $excel.ActiveWorkbook

# --- Employees sheet -------------------------------------------------
$wsEmployees = $wb.Worksheets.Item("Employees")

# min_hours_per_week -> 30 for every current employee (rows 2-61)
$wsEmployees.Range("I2:I61").Value = 30

# Drop the last two employees (Jacob, Leah)
$wsEmployees.Range("A62:O63").EntireRow.Delete() | Out-Null

$wsEmployees.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$wsEmployees.Rows("62:63").Select() | Out-Null

# --- Parameters sheet --------------------------------------------------
$wsParameters = $wb.Worksheets.Item("Parameters")

# minusHours_per_Week -> 0
$wsParameters.Range("B4").Value = 0
$wsParameters.Range("D6").Select() | Out-Null

# Employees ends up the active tab/sheet
$wsEmployees.Activate() | Out-Null
